$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Amount of hours (column B) and Date (column C) for the new hour-log rows 12-16
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 44353

$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 44353

$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 44353

$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 44353

$ws.Range("B16").Value = 1
$ws.Range("C16").Value = 44354

# Subject (column A) and Description (column D) text, entered in the same
# order the author originally typed them so the shared-string table lines
# up with the saved workbook.
$ws.Range("A12").Value = "Router"
$ws.Range("A13").Value = "Eliminate Repetitiveness"
$ws.Range("A15").Value = "Fixing of specific functions"

$ws.Range("D12").Value = "Complete setup of the router for Listings, Houses, Rooms and Users"
$ws.Range("D13").Value = "Replace declaration of Request and Response in every Controller. Controllers inherit from BaseController instead of Mvc/Controller. View is disabled in BaseController before execution, instead of every function in every controller. Added an extra function in Controller Base."
$ws.Range("D14").Value = "Renamed functions and variables that had unproper naming. Removed commented code, corrected comments that had wrong information (e.g. mentioning the user when creating a new house). "

$ws.Range("A14").Value = "Comments/ variables fixing"

$ws.Range("D15").Value = "Replaced 'findFirst(id = ?)' with 'findFirstById(?)'. ControllerBase has the 'errorCheck' function. When the check is passed a generic response message is given."

$ws.Range("A16").Value = "Extra Documentation"

$ws.Range("D16").Value = "Added Router documentation and Postman JSON file. Readme file contains more information."

# Move the active selection as recorded in the saved workbook
$ws.Range("D16").Select() | Out-Null
